$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values for the two middle detail rows (E16 and E18),
# leaving E17 unchanged (1-base visual swap matching the updated account statement).
$ws.Range("E16").Value = "2102"
$ws.Range("E18").Value = "2104"

# Update the "Valor Mora" amount for the last worker row (new period 2109 entry).
$ws.Range("G19").Value = 1423500
